# Updated symbol list on Sun Dec 25 11:45:50 UTC 2022 with GitHub Actions
#
# The "Price" column (D) and a couple of "Volume(1h)" label cells (E) are
# stored as plain text (inline strings) in the workbook, not numbers.
# When writing a numeric-looking string back through Range.Value, Excel's
# normal "smart" input parsing would convert it to a real number, which
# would not match the source data (and would also silently re-type the
# cell). Prefixing the value with a single quote forces it to stay text,
# exactly like typing '244.91 into a cell in Excel. That, however, also
# stamps the cell with a "quote prefix" display flag, so we immediately
# reset the cell's style back to Normal — this keeps the cell as General/
# unstyled (matching the original workbook) while preserving the text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice([string]$cellRef, [string]$newValue) {
    $ws.Range($cellRef).Value = "'" + $newValue
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextPrice "D2"  "244.91"
Set-TextPrice "D3"  "23.06"
Set-TextPrice "D4"  "5.437"
Set-TextPrice "D5"  "0.05968"
Set-TextPrice "D6"  "3.390"
Set-TextPrice "D7"  "0.8100"
Set-TextPrice "D8"  "0.9263"
Set-TextPrice "D10" "0.07437"
Set-TextPrice "D11" "0.03394"
Set-TextPrice "D12" "0.03037"
Set-TextPrice "D13" "0.09339"
Set-TextPrice "D14" "3.962"
Set-TextPrice "D15" "0.001595"
Set-TextPrice "D17" "0.0005942"
Set-TextPrice "D18" "0.005423"
Set-TextPrice "D19" "0.004158"
Set-TextPrice "D20" "0.0009828"
Set-TextPrice "D21" "0.00007703"
Set-TextPrice "D23" "6.453"

Set-TextPrice "D41" "0.006212"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextPrice "D44" "0.007122"
$ws.Range("E44").Value = "43LocalTradersLCT"

Set-TextPrice "D45" "0.00005173"
Set-TextPrice "D47" "0.0005802"
Set-TextPrice "D48" "1.080"
